$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = '@'
$ws.Cells.Item(2, 4).Value = '51.352.55'
$ws.Cells.Item(2, 4).Style = 'Normal'
$ws.Cells.Item(2, 5).NumberFormat = '@'
$ws.Cells.Item(2, 5).Value = '  -0.68%  '
$ws.Cells.Item(2, 5).Style = 'Normal'
$ws.Cells.Item(3, 4).NumberFormat = '@'
$ws.Cells.Item(3, 4).Value = '3.100.30'
$ws.Cells.Item(3, 4).Style = 'Normal'
$ws.Cells.Item(3, 5).NumberFormat = '@'
$ws.Cells.Item(3, 5).Value = '  +2.12%  '
$ws.Cells.Item(3, 5).Style = 'Normal'
$ws.Cells.Item(4, 5).NumberFormat = '@'
$ws.Cells.Item(4, 5).Value = '  -0.07%  '
$ws.Cells.Item(4, 5).Style = 'Normal'
$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '389.02'
$ws.Cells.Item(5, 4).Style = 'Normal'
$ws.Cells.Item(5, 5).NumberFormat = '@'
$ws.Cells.Item(5, 5).Value = '  +2.12%  '
$ws.Cells.Item(5, 5).Style = 'Normal'
$ws.Cells.Item(6, 4).NumberFormat = '@'
$ws.Cells.Item(6, 4).Value = '104.07'
$ws.Cells.Item(6, 4).Style = 'Normal'
$ws.Cells.Item(6, 5).NumberFormat = '@'
$ws.Cells.Item(6, 5).Value = '  +1.12%  '
$ws.Cells.Item(6, 5).Style = 'Normal'
$ws.Cells.Item(7, 4).NumberFormat = '@'
$ws.Cells.Item(7, 4).Value = '0.539'
$ws.Cells.Item(7, 4).Style = 'Normal'
$ws.Cells.Item(7, 5).NumberFormat = '@'
$ws.Cells.Item(7, 5).Value = '  -1.70%  '
$ws.Cells.Item(7, 5).Style = 'Normal'
$ws.Cells.Item(8, 5).NumberFormat = '@'
$ws.Cells.Item(8, 5).Value = '  -0.01%  '
$ws.Cells.Item(8, 5).Style = 'Normal'
$ws.Cells.Item(9, 5).NumberFormat = '@'
$ws.Cells.Item(9, 5).Value = '  -1.10%  '
$ws.Cells.Item(9, 5).Style = 'Normal'
$ws.Cells.Item(10, 4).NumberFormat = '@'
$ws.Cells.Item(10, 4).Value = '36.94'
$ws.Cells.Item(10, 4).Style = 'Normal'
$ws.Cells.Item(10, 5).NumberFormat = '@'
$ws.Cells.Item(10, 5).Value = '  +0.32%  '
$ws.Cells.Item(10, 5).Style = 'Normal'
$ws.Cells.Item(11, 5).NumberFormat = '@'
$ws.Cells.Item(11, 5).Value = '  +0.05%  '
$ws.Cells.Item(11, 5).Style = 'Normal'
$ws.Cells.Item(12, 5).NumberFormat = '@'
$ws.Cells.Item(12, 5).Value = '  -1.02%  '
$ws.Cells.Item(12, 5).Style = 'Normal'
$ws.Cells.Item(13, 4).NumberFormat = '@'
$ws.Cells.Item(13, 4).Value = '3.587.30'
$ws.Cells.Item(13, 4).Style = 'Normal'
$ws.Cells.Item(13, 5).NumberFormat = '@'
$ws.Cells.Item(13, 5).Value = '  +2.13%  '
$ws.Cells.Item(13, 5).Style = 'Normal'
$ws.Cells.Item(14, 4).NumberFormat = '@'
$ws.Cells.Item(14, 4).Value = '18.54'
$ws.Cells.Item(14, 4).Style = 'Normal'
$ws.Cells.Item(14, 5).NumberFormat = '@'
$ws.Cells.Item(14, 5).Value = '  +0.11%  '
$ws.Cells.Item(14, 5).Style = 'Normal'
$ws.Cells.Item(15, 4).NumberFormat = '@'
$ws.Cells.Item(15, 4).Value = '7.75'
$ws.Cells.Item(15, 4).Style = 'Normal'
$ws.Cells.Item(15, 5).NumberFormat = '@'
$ws.Cells.Item(15, 5).Value = '  -0.05%  '
$ws.Cells.Item(15, 5).Style = 'Normal'
$ws.Cells.Item(16, 4).NumberFormat = '@'
$ws.Cells.Item(16, 4).Value = '3.106.91'
$ws.Cells.Item(16, 4).Style = 'Normal'
$ws.Cells.Item(16, 5).NumberFormat = '@'
$ws.Cells.Item(16, 5).Value = '  +2.66%  '
$ws.Cells.Item(16, 5).Style = 'Normal'
$ws.Cells.Item(17, 5).NumberFormat = '@'
$ws.Cells.Item(17, 5).Value = '  +1.88%  '
$ws.Cells.Item(17, 5).Style = 'Normal'
$ws.Cells.Item(18, 4).NumberFormat = '@'
$ws.Cells.Item(18, 4).Value = '10.72'
$ws.Cells.Item(18, 4).Style = 'Normal'
$ws.Cells.Item(18, 5).NumberFormat = '@'
$ws.Cells.Item(18, 5).Value = '  +2.13%  '
$ws.Cells.Item(18, 5).Style = 'Normal'
$ws.Cells.Item(19, 4).NumberFormat = '@'
$ws.Cells.Item(19, 4).Value = '51.444.66'
$ws.Cells.Item(19, 4).Style = 'Normal'
$ws.Cells.Item(19, 5).NumberFormat = '@'
$ws.Cells.Item(19, 5).Value = '  -0.58%  '
$ws.Cells.Item(19, 5).Style = 'Normal'
$ws.Cells.Item(20, 5).NumberFormat = '@'
$ws.Cells.Item(20, 5).Value = '  +4.61%  '
$ws.Cells.Item(20, 5).Style = 'Normal'
$ws.Cells.Item(21, 4).NumberFormat = '@'
$ws.Cells.Item(21, 4).Value = '12.44'
$ws.Cells.Item(21, 4).Style = 'Normal'
$ws.Cells.Item(21, 5).NumberFormat = '@'
$ws.Cells.Item(21, 5).Value = '  -0.30%  '
$ws.Cells.Item(21, 5).Style = 'Normal'
$ws.Cells.Item(22, 5).NumberFormat = '@'
$ws.Cells.Item(22, 5).Value = '  +0.03%  '
$ws.Cells.Item(22, 5).Style = 'Normal'
$ws.Cells.Item(23, 4).NumberFormat = '@'
$ws.Cells.Item(23, 4).Value = '70.01'
$ws.Cells.Item(23, 4).Style = 'Normal'
$ws.Cells.Item(23, 5).NumberFormat = '@'
$ws.Cells.Item(23, 5).Value = '  -0.11%  '
$ws.Cells.Item(23, 5).Style = 'Normal'
$ws.Cells.Item(24, 4).NumberFormat = '@'
$ws.Cells.Item(24, 4).Value = '266.01'
$ws.Cells.Item(24, 4).Style = 'Normal'
$ws.Cells.Item(24, 5).NumberFormat = '@'
$ws.Cells.Item(24, 5).Value = '  -0.66%  '
$ws.Cells.Item(24, 5).Style = 'Normal'
$ws.Cells.Item(25, 4).NumberFormat = '@'
$ws.Cells.Item(25, 4).Value = '3.17'
$ws.Cells.Item(25, 4).Style = 'Normal'
$ws.Cells.Item(25, 5).NumberFormat = '@'
$ws.Cells.Item(25, 5).Value = '  +0.15%  '
$ws.Cells.Item(25, 5).Style = 'Normal'
$ws.Cells.Item(26, 4).NumberFormat = '@'
$ws.Cells.Item(26, 4).Value = '7.99'
$ws.Cells.Item(26, 4).Style = 'Normal'
$ws.Cells.Item(26, 5).NumberFormat = '@'
$ws.Cells.Item(26, 5).Value = '  -3.43%  '
$ws.Cells.Item(26, 5).Style = 'Normal'
$ws.Cells.Item(27, 4).NumberFormat = '@'
$ws.Cells.Item(27, 4).Value = '27.42'
$ws.Cells.Item(27, 4).Style = 'Normal'
$ws.Cells.Item(27, 5).NumberFormat = '@'
$ws.Cells.Item(27, 5).Value = '  +4.42%  '
$ws.Cells.Item(27, 5).Style = 'Normal'
$ws.Cells.Item(28, 4).NumberFormat = '@'
$ws.Cells.Item(28, 4).Value = '7.25'
$ws.Cells.Item(28, 4).Style = 'Normal'
$ws.Cells.Item(28, 5).NumberFormat = '@'
$ws.Cells.Item(28, 5).Value = '  -5.51%  '
$ws.Cells.Item(28, 5).Style = 'Normal'
$ws.Cells.Item(29, 5).NumberFormat = '@'
$ws.Cells.Item(29, 5).Value = '  +0.06%  '
$ws.Cells.Item(29, 5).Style = 'Normal'
$ws.Cells.Item(30, 5).NumberFormat = '@'
$ws.Cells.Item(30, 5).Value = '  -3.83%  '
$ws.Cells.Item(30, 5).Style = 'Normal'
$ws.Cells.Item(31, 5).NumberFormat = '@'
$ws.Cells.Item(31, 5).Value = '  -1.57%  '
$ws.Cells.Item(31, 5).Style = 'Normal'
$ws.Cells.Item(32, 4).NumberFormat = '@'
$ws.Cells.Item(32, 4).Value = '10.42'
$ws.Cells.Item(32, 4).Style = 'Normal'
$ws.Cells.Item(32, 5).NumberFormat = '@'
$ws.Cells.Item(32, 5).Value = '  +1.27%  '
$ws.Cells.Item(32, 5).Style = 'Normal'
$ws.Cells.Item(33, 4).NumberFormat = '@'
$ws.Cells.Item(33, 4).Value = '36.25'
$ws.Cells.Item(33, 4).Style = 'Normal'
$ws.Cells.Item(33, 5).NumberFormat = '@'
$ws.Cells.Item(33, 5).Value = '  +7.06%  '
$ws.Cells.Item(33, 5).Style = 'Normal'
$ws.Cells.Item(34, 4).NumberFormat = '@'
$ws.Cells.Item(34, 4).Value = '0.0476'
$ws.Cells.Item(34, 4).Style = 'Normal'
$ws.Cells.Item(34, 5).NumberFormat = '@'
$ws.Cells.Item(34, 5).Value = '  +6.46%  '
$ws.Cells.Item(34, 5).Style = 'Normal'
$ws.Cells.Item(35, 5).NumberFormat = '@'
$ws.Cells.Item(35, 5).Value = '  -0.45%  '
$ws.Cells.Item(35, 5).Style = 'Normal'
$ws.Cells.Item(36, 4).NumberFormat = '@'
$ws.Cells.Item(36, 4).Value = '49.82'
$ws.Cells.Item(36, 4).Style = 'Normal'
$ws.Cells.Item(36, 5).NumberFormat = '@'
$ws.Cells.Item(36, 5).Value = '  -1.43%  '
$ws.Cells.Item(36, 5).Style = 'Normal'
$ws.Cells.Item(37, 5).NumberFormat = '@'
$ws.Cells.Item(37, 5).Value = '  -0.09%  '
$ws.Cells.Item(37, 5).Style = 'Normal'
$ws.Cells.Item(38, 4).NumberFormat = '@'
$ws.Cells.Item(38, 4).Value = '3.40'
$ws.Cells.Item(38, 4).Style = 'Normal'
$ws.Cells.Item(38, 5).NumberFormat = '@'
$ws.Cells.Item(38, 5).Value = '  +2.41%  '
$ws.Cells.Item(38, 5).Style = 'Normal'
$ws.Cells.Item(39, 4).NumberFormat = '@'
$ws.Cells.Item(39, 4).Value = '0.291'
$ws.Cells.Item(39, 4).Style = 'Normal'
$ws.Cells.Item(39, 5).NumberFormat = '@'
$ws.Cells.Item(39, 5).Value = '  -2.85%  '
$ws.Cells.Item(39, 5).Style = 'Normal'
$ws.Cells.Item(40, 4).NumberFormat = '@'
$ws.Cells.Item(40, 4).Value = '130.58'
$ws.Cells.Item(40, 4).Style = 'Normal'
$ws.Cells.Item(40, 5).NumberFormat = '@'
$ws.Cells.Item(40, 5).Value = '  +2.21%  '
$ws.Cells.Item(40, 5).Style = 'Normal'
$ws.Cells.Item(41, 5).NumberFormat = '@'
$ws.Cells.Item(41, 5).Value = '  -0.57%  '
$ws.Cells.Item(41, 5).Style = 'Normal'
$ws.Cells.Item(42, 4).NumberFormat = '@'
$ws.Cells.Item(42, 4).Value = '3.88'
$ws.Cells.Item(42, 4).Style = 'Normal'
$ws.Cells.Item(42, 5).NumberFormat = '@'
$ws.Cells.Item(42, 5).Value = '  +2.69%  '
$ws.Cells.Item(42, 5).Style = 'Normal'
$ws.Cells.Item(43, 4).NumberFormat = '@'
$ws.Cells.Item(43, 4).Value = '16.66'
$ws.Cells.Item(43, 4).Style = 'Normal'
$ws.Cells.Item(43, 5).NumberFormat = '@'
$ws.Cells.Item(43, 5).Value = '  -2.18%  '
$ws.Cells.Item(43, 5).Style = 'Normal'
$ws.Cells.Item(44, 5).NumberFormat = '@'
$ws.Cells.Item(44, 5).Value = '  -0.29%  '
$ws.Cells.Item(44, 5).Style = 'Normal'
$ws.Cells.Item(45, 5).NumberFormat = '@'
$ws.Cells.Item(45, 5).Value = '  -2.06%  '
$ws.Cells.Item(45, 5).Style = 'Normal'
$ws.Cells.Item(46, 4).NumberFormat = '@'
$ws.Cells.Item(46, 4).Value = '22.10'
$ws.Cells.Item(46, 4).Style = 'Normal'
$ws.Cells.Item(46, 5).NumberFormat = '@'
$ws.Cells.Item(46, 5).Value = '  +1.74%  '
$ws.Cells.Item(46, 5).Style = 'Normal'
$ws.Cells.Item(47, 4).NumberFormat = '@'
$ws.Cells.Item(47, 4).Value = '2.51'
$ws.Cells.Item(47, 4).Style = 'Normal'
$ws.Cells.Item(47, 5).NumberFormat = '@'
$ws.Cells.Item(47, 5).Value = '  +4.60%  '
$ws.Cells.Item(47, 5).Style = 'Normal'
$ws.Cells.Item(48, 5).NumberFormat = '@'
$ws.Cells.Item(48, 5).Value = '  -0.30%  '
$ws.Cells.Item(48, 5).Style = 'Normal'
$ws.Cells.Item(49, 4).NumberFormat = '@'
$ws.Cells.Item(49, 4).Value = '2.073.87'
$ws.Cells.Item(49, 4).Style = 'Normal'
$ws.Cells.Item(49, 5).NumberFormat = '@'
$ws.Cells.Item(49, 5).Value = '  +2.09%  '
$ws.Cells.Item(49, 5).Style = 'Normal'
$ws.Cells.Item(50, 4).NumberFormat = '@'
$ws.Cells.Item(50, 4).Value = '0.937'
$ws.Cells.Item(50, 4).Style = 'Normal'
$ws.Cells.Item(50, 5).NumberFormat = '@'
$ws.Cells.Item(50, 5).Value = '  +19.18%  '
$ws.Cells.Item(50, 5).Style = 'Normal'
$ws.Cells.Item(51, 4).NumberFormat = '@'
$ws.Cells.Item(51, 4).Value = '0.0327'
$ws.Cells.Item(51, 4).Style = 'Normal'
$ws.Cells.Item(51, 5).NumberFormat = '@'
$ws.Cells.Item(51, 5).Value = '  +1.59%  '
$ws.Cells.Item(51, 5).Style = 'Normal'
